# Add new columns I ("I0") and J ("IF") to the sheet, matching the
# existing header style used in column H, and fill in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell H1 onto the two new header cells, then set their
# text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Numeric data for rows 2-45 (column I then column J).
$iValues = @(7,4,7,3,6,7,7,5,9,7,8,8,5,11,7,7,7,6,8,4,8,8,7,8,12,6,7,6,8,8,7,8,6,7,8,6,9,9,9,9,6,5,5,5)
$jValues = @(8,5,7,4,7,7,8,6,9,8,8,8,7,11,8,7,8,7,8,6,8,8,7,8,13,6,8,6,8,8,7,8,7,8,8,7,9,9,9,9,6,5,5,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Host "Added columns I (I0) and J (IF) with data for rows 1-45"
